# daily auto push: 2026-01-10 13:37 UTC
# Insert two new observation rows for 2026/01/10 right after the existing
# 2026/01/10 block (which ended at row 594), pushing all subsequent rows
# down by two. This grows the sheet from A1:D636 to A1:D638.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 595-596; everything from the old row 595 onward
# (2026/12/29 ... 2027/01/05) shifts down to rows 597-638.
$ws.Rows("595:596").Insert()

# Format the date column as Text first so the "yyyy/mm/dd"-looking string
# is stored literally (matching every other row in column A) instead of
# being auto-converted into a date serial number by Excel's input parser.
$ws.Range("A595:A596").NumberFormat = "@"

$ws.Range("A595").Value = "2026/01/10"
$ws.Range("B595").Value = "土"
$ws.Range("C595").Value = 18
$ws.Range("D595").Value = 32

$ws.Range("A596").Value = "2026/01/10"
$ws.Range("B596").Value = "土"
$ws.Range("C596").Value = 20
$ws.Range("D596").Value = 33

# Restore the default (unstyled) look so the new rows don't carry a
# leftover "Text" number format, consistent with the rest of the sheet.
$ws.Range("A595:A596").Style = "Normal"

Write-Host "Inserted rows 595-596 (2026/01/10); sheet now spans to row 638."
